$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (dates 44432-44440, i.e. 2021-08-24 .. 2021-09-01)
$data = @(
    @(44432, 0, 5, 94.6969696969697),
    @(44433, 1, 6, 113.6363636363636),
    @(44434, 0, 4, 75.75757575757575),
    @(44435, 1, 4, 75.75757575757575),
    @(44436, 1, 4, 75.75757575757575),
    @(44437, 0, 4, 75.75757575757575),
    @(44438, 1, 4, 75.75757575757575),
    @(44439, 0, 4, 75.75757575757575),
    @(44440, 0, 3, 56.81818181818181)
)

$startRow = 358
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]

    $cellA = $ws.Cells.Item($row, 1)
    $ws.Cells.Item($row - 1, 1).Copy() | Out-Null
    $cellA.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $cellA.Value = $vals[0]

    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
}
